$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value, per the target diff
$changes = @{
    "F6"  = 1
    "F8"  = -1
    "F9"  = -2
    "E13" = 1
    "F13" = -2
    "H13" = 3
    "I13" = 7
    "F19" = 2
    "F23" = -2
    "F32" = -4
    "F34" = 1
    "F38" = 3
    "F40" = -2
    "F41" = 1
    "F42" = -2
    "F44" = -2
    "F51" = -3
    "F52" = -2
    "F58" = -4
    "F61" = -9
    "F64" = 1
    "F65" = -6
    "F69" = -6
    "F73" = 2
    "F74" = -1
    "F75" = 1
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
